$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range('D2').Value = '68.535.58'
$ws.Range('E2').Value = '  -1.51%  '
$ws.Range('D3').Value = '3.858.80'
$ws.Range('E3').Value = '  -0.74%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '602.68'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.07'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.15%  '
$ws.Range('D7').Value = '3.857.74'
$ws.Range('E7').Value = '  -0.74%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -1.09%  '
$ws.Range('E10').Value = '  -1.67%  '
$ws.Range('E11').Value = '  +1.14%  '
$ws.Range('E12').Value = '  -2.23%  '
$ws.Range('E13').Value = '  +4.29%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.13'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.85%  '
$ws.Range('D15').Value = '4.505.12'
$ws.Range('E15').Value = '  -0.86%  '
$ws.Range('D16').Value = '3.848.75'
$ws.Range('E16').Value = '  -1.64%  '
$ws.Range('D17').Value = '68.692.24'
$ws.Range('E17').Value = '  -1.35%  '
$ws.Range('E18').Value = '  -0.76%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.39'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.95%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.45'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.66%  '
$ws.Range('E21').Value = '  -0.91%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '471.27'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.82%  '
$ws.Range('E23').Value = '  -1.41%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.58'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.08%  '
$ws.Range('E26').Value = '  -2.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.14'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.53%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.24'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.05%  '
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('E30').Value = '  -0.61%  '
$ws.Range('D31').Value = '4.009.38'
$ws.Range('E31').Value = '  -0.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.69'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.86%  '
$ws.Range('E33').Value = '  -1.14%  '
$ws.Range('E34').Value = '  -3.67%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.35'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.85%  '
$ws.Range('D36').Value = '3.823.33'
$ws.Range('E36').Value = '  -0.85%  '
$ws.Range('E37').Value = '  -2.35%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.72'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +8.99%  '
$ws.Range('E39').Value = '  -1.51%  '
$ws.Range('E40').Value = '  -2.65%  '
$ws.Range('E41').Value = '  -2.87%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('E43').Value = '  -3.69%  '
$ws.Range('E44').Value = '  -4.68%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.75'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.69%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '416.90'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.23%  '
$ws.Range('B47').Value = 'FLOKI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.000294'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +6.27%  '
$ws.Range('E48').Value = '  +0.01%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '47.12'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.01%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '142.10'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.12%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0360'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.58%  '
